# Added ifoCAST full series evaluation: for each data row, the error
# series is shifted one column to the left (the old leading value -
# matched to the quarter that dropped out of the ifoCAST window - is
# discarded) and a freshly evaluated error is appended at the end of the
# still-full rows; rows that have no new observation to append simply
# lose their trailing column, shrinking the triangular matrix by one
# cell per row as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2177158050940863
$ws.Range("C2").Value = -0.1220406610940863
$ws.Range("D2").Value = -0.3824598050940862
$ws.Range("E2").Value = 0.4529951949059138
$ws.Range("F2").Value = 0.3529222419059138
$ws.Range("G2").Value = -0.0008708050940862688
$ws.Range("H2").Value = 0.1473071949059137
$ws.Range("I2").Value = 0.3047531949059137
$ws.Range("J2").Value = -0.1959098050940863
$ws.Range("K2").Value = 0.08445119490591374
$ws.Range("B3").Value = -0.1395947824540507
$ws.Range("C3").Value = -0.4000139264540507
$ws.Range("D3").Value = 0.4354410735459493
$ws.Range("E3").Value = 0.3353681205459493
$ws.Range("F3").Value = -0.01842492645405072
$ws.Range("G3").Value = 0.1297530735459493
$ws.Range("H3").Value = 0.2871990735459493
$ws.Range("I3").Value = -0.2134639264540507
$ws.Range("J3").Value = 0.06689707354594927
$ws.Range("K3").Value = -0.1050779264540507
$ws.Range("B4").Value = -0.3119066075030444
$ws.Range("C4").Value = 0.5235483924969555
$ws.Range("D4").Value = 0.4234754394969555
$ws.Range("E4").Value = 0.06968239249695551
$ws.Range("F4").Value = 0.2178603924969555
$ws.Range("G4").Value = 0.3753063924969555
$ws.Range("H4").Value = -0.1253566075030445
$ws.Range("I4").Value = 0.1550043924969555
$ws.Range("J4").Value = -0.01697060750304449
$ws.Range("K4").Value = 0.2888923924969555
$ws.Range("B5").Value = 0.7021230259847391
$ws.Range("C5").Value = 0.6020500729847391
$ws.Range("D5").Value = 0.2482570259847391
$ws.Range("E5").Value = 0.3964350259847391
$ws.Range("F5").Value = 0.5538810259847391
$ws.Range("G5").Value = 0.0532180259847391
$ws.Range("H5").Value = 0.3335790259847391
$ws.Range("I5").Value = 0.1616040259847391
$ws.Range("J5").Value = 0.4674670259847391
$ws.Range("K5").Value = -0.1448629740152609
$ws.Range("B6").Value = 1.514077300737389
$ws.Range("C6").Value = 1.160284253737389
$ws.Range("D6").Value = 1.308462253737389
$ws.Range("E6").Value = 1.465908253737389
$ws.Range("F6").Value = 0.9652452537373891
$ws.Range("G6").Value = 1.245606253737389
$ws.Range("H6").Value = 1.073631253737389
$ws.Range("I6").Value = 1.379494253737389
$ws.Range("J6").Value = 0.7671642537373891
$ws.Range("K6").Value = 1.455535253737389
$ws.Range("B7").Value = 0.2163100177716323
$ws.Range("C7").Value = 0.3644880177716323
$ws.Range("D7").Value = 0.5219340177716323
$ws.Range("E7").Value = 0.02127101777163229
$ws.Range("F7").Value = 0.3016320177716323
$ws.Range("G7").Value = 0.1296570177716323
$ws.Range("H7").Value = 0.4355200177716323
$ws.Range("I7").Value = -0.1768099822283677
$ws.Range("J7").Value = 0.5115610177716323
$ws.Range("K7").ClearContents()
$ws.Range("B8").Value = 0.3684559152847414
$ws.Range("C8").Value = 0.5259019152847414
$ws.Range("D8").Value = 0.0252389152847414
$ws.Range("E8").Value = 0.3055999152847414
$ws.Range("F8").Value = 0.1336249152847414
$ws.Range("G8").Value = 0.4394879152847414
$ws.Range("H8").Value = -0.1728420847152586
$ws.Range("I8").Value = 0.5155289152847414
$ws.Range("J8").ClearContents()
$ws.Range("B9").Value = 0.6615420054549828
$ws.Range("C9").Value = 0.1608790054549828
$ws.Range("D9").Value = 0.4412400054549828
$ws.Range("E9").Value = 0.2692650054549828
$ws.Range("F9").Value = 0.5751280054549828
$ws.Range("G9").Value = -0.03720199454501721
$ws.Range("H9").Value = 0.6511690054549828
$ws.Range("I9").ClearContents()
$ws.Range("B10").Value = -0.07992717640068192
$ws.Range("C10").Value = 0.2004338235993181
$ws.Range("D10").Value = 0.02845882359931809
$ws.Range("E10").Value = 0.3343218235993181
$ws.Range("F10").Value = -0.2780081764006819
$ws.Range("G10").Value = 0.4103628235993181
$ws.Range("H10").ClearContents()
$ws.Range("B11").Value = 0.1551026083472217
$ws.Range("C11").Value = -0.01687239165277828
$ws.Range("D11").Value = 0.2889906083472217
$ws.Range("E11").Value = -0.3233393916527783
$ws.Range("F11").Value = 0.3650316083472218
$ws.Range("G11").ClearContents()
$ws.Range("B12").Value = -0.08373381358883215
$ws.Range("C12").Value = 0.2221291864111679
$ws.Range("D12").Value = -0.3902008135888321
$ws.Range("E12").Value = 0.2981701864111679
$ws.Range("F12").ClearContents()
$ws.Range("B13").Value = 0.1925429840544512
$ws.Range("C13").Value = -0.4197870159455487
$ws.Range("D13").Value = 0.2685839840544513
$ws.Range("E13").ClearContents()
$ws.Range("B14").Value = -0.4379376726751896
$ws.Range("C14").Value = 0.2504333273248104
$ws.Range("D14").ClearContents()
$ws.Range("B15").Value = 0.2324015030761111
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
